$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws 'D2' '292.73'
Set-TextValue $ws 'E2' '-7.05%'
Set-TextValue $ws 'E3' '-1.19%'
Set-TextValue $ws 'D4' '5.043'
Set-TextValue $ws 'E4' '-2.12%'
Set-TextValue $ws 'D5' '0.07325'
Set-TextValue $ws 'E5' '-3.70%'
Set-TextValue $ws 'D6' '1.535'
Set-TextValue $ws 'E6' '-9.86%'
Set-TextValue $ws 'D7' '0.9303'
Set-TextValue $ws 'D8' '2.369'
Set-TextValue $ws 'E8' '-2.27%'
Set-TextValue $ws 'D9' '0.1170'
Set-TextValue $ws 'E9' '-2.58%'
Set-TextValue $ws 'D10' '0.1742'
Set-TextValue $ws 'D11' '0.04336'
Set-TextValue $ws 'E11' '4.77%'
Set-TextValue $ws 'D12' '0.08672'
Set-TextValue $ws 'E12' '-3.62%'
Set-TextValue $ws 'D13' '0.1053'
Set-TextValue $ws 'E13' '0.09%'
Set-TextValue $ws 'B14' 'TigerCash'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D14' '0.006001'
Set-TextValue $ws 'E14' '2.71%'
Set-TextValue $ws 'B15' 'LEO'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D15' '3.337'
Set-TextValue $ws 'E15' '0.04%'
Set-TextValue $ws 'B16' 'GateToken'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D16' '4.284'
Set-TextValue $ws 'E16' '-0.94%'
Set-TextValue $ws 'B17' 'BitpandaEcosystemToken'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws 'D17' '0.3289'
Set-TextValue $ws 'E17' '-1.79%'
Set-TextValue $ws 'B18' 'MCDex'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws 'D18' '7.970'
Set-TextValue $ws 'E18' '5.05%'
Set-TextValue $ws 'B19' 'ProBitToken'
Set-TextValue $ws 'C19' 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws 'D19' '0.1400'
Set-TextValue $ws 'E19' '3.66%'
Set-TextValue $ws 'B20' 'ZBToken'
Set-TextValue $ws 'C20' 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws 'D20' '0.2742'
Set-TextValue $ws 'E20' '-3.35%'
Set-TextValue $ws 'B21' 'BitForexToken'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D21' '0.001266'
Set-TextValue $ws 'E21' '-1.65%'
Set-TextValue $ws 'D22' '0.03941'
Set-TextValue $ws 'E22' '-0.58%'
Set-TextValue $ws 'E23' '-1.50%'
Set-TextValue $ws 'D24' '0.003789'
Set-TextValue $ws 'E24' '-7.14%'
Set-TextValue $ws 'E25' '-4.99%'
Set-TextValue $ws 'D26' '0.0003724'
Set-TextValue $ws 'D38' '0.02302'
Set-TextValue $ws 'E38' '-4.85%'
Set-TextValue $ws 'D39' '0.05064'
Set-TextValue $ws 'E39' '-2.29%'
Set-TextValue $ws 'D40' '0.006224'
Set-TextValue $ws 'E40' '88.60%'
Set-TextValue $ws 'D41' '0.007819'
Set-TextValue $ws 'E41' '1.47%'
Set-TextValue $ws 'D42' '0.1288'
Set-TextValue $ws 'E42' '-1.15%'
Set-TextValue $ws 'D43' '0.007367'
Set-TextValue $ws 'E43' '-2.90%'
Set-TextValue $ws 'D44' '0.008213'
Set-TextValue $ws 'E44' '-3.09%'
Set-TextValue $ws 'D45' '0.2925'
Set-TextValue $ws 'E45' '-13.77%'
Set-TextValue $ws 'D46' '0.00006281'
Set-TextValue $ws 'E46' '-4.76%'
Set-TextValue $ws 'D47' '0.00000000750'
Set-TextValue $ws 'E47' '0.06%'
Set-TextValue $ws 'D48' '0.03375'
Set-TextValue $ws 'E48' '-87.44%'
Set-TextValue $ws 'D49' '0.00002101'
Set-TextValue $ws 'E49' '0.06%'
Set-TextValue $ws 'E50' '0.06%'
